# Fruta / hortaliza, semanal
# Insert two new weekly-report rows for "Femacal de La Calera" / Frutilla
# right before the existing row 252, shifting the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 252 (old rows 252-269 become 254-271)
$ws.Rows("252:253").Insert()

# Constant columns shared by every data row in this block
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "`$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7

function Set-DataRow($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $tipo
    $ws.Cells.Item($Row, 7).Value = $productoId
    $ws.Cells.Item($Row, 8).Value = $producto
    $ws.Cells.Item($Row, 9).Value = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

# New row 252: "Especial" quality, 2022-02-18
Set-DataRow 252 44610 "Especial" 56 7000 7000 7000 1000

# New row 253: "Segunda" quality, 2022-02-18
Set-DataRow 253 44610 "Segunda" 40 4000 4000 4000 571
